# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) across several Leve-profit sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 897
$ws.Range("I6").Value = 1137.2727
$ws.Range("K6").Value = 3411.8181
$ws.Range("M6").Value = -3299.8181

$ws.Range("H40").Value = 7340.4116
$ws.Range("J40").Value = 8071.5454
$ws.Range("L40").Value = 8071.5454
$ws.Range("N40").Value = -8421.545399999999

$ws.Range("H107").Value = 48178.477
$ws.Range("I107").Value = 48178.477
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 48178.477
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -46258.477
$ws.Range("N107").ClearContents()

$ws.Range("H112").Value = 2123.2173
$ws.Range("J112").Value = 2123.2173
$ws.Range("L112").Value = 6369.651899999999
$ws.Range("N112").Value = -8585.651899999999

$ws.Range("H113").Value = 2100
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254

$ws.Range("H129").Value = 2122.2917
$ws.Range("I129").Value = 540.375
$ws.Range("J129").Value = 2913.25
$ws.Range("K129").Value = 1621.125
$ws.Range("L129").Value = 8739.75
$ws.Range("M129").Value = 3378.875
$ws.Range("N129").Value = -18739.75

$ws.Range("H132").Value = 1747.4857
$ws.Range("I132").Value = 1938.7241
$ws.Range("J132").Value = 823.1667
$ws.Range("K132").Value = 5816.1723
$ws.Range("L132").Value = 2469.5001
$ws.Range("M132").Value = -3286.1723
$ws.Range("N132").Value = -7529.5001

$ws.Range("H135").Value = 456316.47
$ws.Range("I135").Value = 668295.5600000001
$ws.Range("J135").Value = 2075.5715
$ws.Range("K135").Value = 6014660.040000001
$ws.Range("L135").Value = 18680.1435
$ws.Range("M135").Value = -6012125.040000001
$ws.Range("N135").Value = -23750.1435

$ws.Range("H138").Value = 4349
$ws.Range("I138").Value = 1545.4736
$ws.Range("J138").Value = 6013.5938
$ws.Range("K138").Value = 4636.4208
$ws.Range("L138").Value = 18040.7814
$ws.Range("M138").Value = 503.5792000000001
$ws.Range("N138").Value = -28320.7814

$ws.Range("H141").Value = 3279.1667
$ws.Range("I141").Value = 3279.1667
$ws.Range("K141").Value = 9837.500100000001
$ws.Range("M141").Value = -4657.500100000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4392.754
$ws.Range("I32").Value = 3991.5615
$ws.Range("K32").Value = 3991.5615
$ws.Range("M32").Value = -3704.5615

$ws.Range("H45").Value = 3257.4167
$ws.Range("I45").Value = 2227.2856
$ws.Range("K45").Value = 2227.2856
$ws.Range("M45").Value = -1850.2856

$ws.Range("H74").Value = 1577.7778
$ws.Range("I74").Value = 1577.7778
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1577.7778
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -703.7778000000001
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1577.7778
$ws.Range("I77").Value = 1577.7778
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7888.889
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3520.889
$ws.Range("N77").ClearContents()

$ws.Range("H102").Value = 1695.4546
$ws.Range("I102").Value = 1894.4445
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 1894.4445
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = -272.4445000000001
$ws.Range("N102").Value = -4044

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1624.8889
$ws.Range("I105").Value = 1589.125
$ws.Range("K105").Value = 1589.125
$ws.Range("M105").Value = 157.875

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3948.4443
$ws.Range("I122").Value = 3961.125
$ws.Range("J122").Value = 3938.3
$ws.Range("K122").Value = 11883.375
$ws.Range("L122").Value = 11814.9
$ws.Range("M122").Value = -9433.375
$ws.Range("N122").Value = -16714.9

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1039.8
$ws.Range("J51").Value = 1039.8
$ws.Range("L51").Value = 3119.4
$ws.Range("N51").Value = -4039.4

$ws.Range("H57").Value = 1466.3334
$ws.Range("I57").Value = 1399
$ws.Range("K57").Value = 4197
$ws.Range("M57").Value = -3638

$ws.Range("H75").Value = 2262.5

$ws.Range("H78").Value = 2262.5

$ws.Range("H94").Value = 8950
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 8950
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 26850
$ws.Range("N94").Value = -28202
$ws.Range("M94").ClearContents()

$ws.Range("H131").Value = 3139.8306
$ws.Range("I131").Value = 1757
$ws.Range("J131").Value = 3240.4
$ws.Range("K131").Value = 5271
$ws.Range("L131").Value = 9721.200000000001
$ws.Range("M131").Value = -231
$ws.Range("N131").Value = -19801.2

$ws.Range("H140").Value = 1655.0333
$ws.Range("I140").Value = 1273.4286
$ws.Range("K140").Value = 3820.2858
$ws.Range("M140").Value = 1359.7142

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 271.625
$ws.Range("J17").Value = 295
$ws.Range("L17").Value = 295
$ws.Range("N17").Value = -631

$ws.Range("H34").Value = 16583.334
$ws.Range("J34").Value = 16583.334
$ws.Range("L34").Value = 16583.334
$ws.Range("N34").Value = -17119.334

$ws.Range("H43").Value = 25998.334
$ws.Range("I43").Value = 13997.5
$ws.Range("K43").Value = 13997.5
$ws.Range("M43").Value = -13846.5

$ws.Range("H76").Value = 16583.334
$ws.Range("J76").Value = 16583.334
$ws.Range("L76").Value = 16583.334
$ws.Range("N76").Value = -17213.334

$ws.Range("H79").Value = 16583.334
$ws.Range("J79").Value = 16583.334
$ws.Range("L79").Value = 16583.334
$ws.Range("N79").Value = -18767.334

$ws.Range("H97").Value = 1553
$ws.Range("I97").Value = 1139.8572
$ws.Range("K97").Value = 1139.8572
$ws.Range("M97").Value = -643.8571999999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H22").Value = 475
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -1190

$ws.Range("H27").Value = 475
$ws.Range("J27").Value = 600
$ws.Range("L27").Value = 600
$ws.Range("N27").Value = -814

$ws.Range("H41").Value = 35000
$ws.Range("J41").Value = 35000
$ws.Range("L41").Value = 35000
$ws.Range("N41").Value = -35876

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2503
$ws.Range("I9").Value = 2503
$ws.Range("K9").Value = 2503
$ws.Range("M9").Value = -2363

$ws.Range("H61").Value = 19750
$ws.Range("I61").Value = 19750
$ws.Range("K61").Value = 19750
$ws.Range("M61").Value = -19458
